# Update the "Pais" (countries) COVID-19 dashboard sheet with the latest
# refresh: new case/death/recovered counts for several countries, and the
# "last updated" timestamp. Two country pairs (Malasia/Guinea and
# Montserrat/Islas Malvinas) swap rank because their updated totals changed
# their relative ordering, so both the label and the data in those rows move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "last updated" timestamp -------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 23 de Septiembre de 2020 a las 12:08"

# --- Row 18: Banglades ------------------------------------------------------
$ws.Range("B18").Value = 353844
$ws.Range("C18").Value = 1666
$ws.Range("D18").Value = 262953
$ws.Range("E18").Value = 85847
$ws.Range("G18").Value = 37
$ws.Range("H18").Value = 5044

# --- Row 24: Filipinas -------------------------------------------------------
$ws.Range("B24").Value = 294591
$ws.Range("C24").Value = 2833
$ws.Range("D24").Value = 231373
$ws.Range("E24").Value = 58127
$ws.Range("G24").Value = 44
$ws.Range("H24").Value = 5091

# --- Row 25: Alemania --------------------------------------------------------
$ws.Range("B25").Value = 277376
$ws.Range("C25").Value = 200
$ws.Range("E25").Value = 19983
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 9493

# --- Row 33: Rumania ---------------------------------------------------------
$ws.Range("B33").Value = 116415
$ws.Range("C33").Value = 1767
$ws.Range("D33").Value = 93558
$ws.Range("E33").Value = 18307
$ws.Range("G33").Value = 47
$ws.Range("H33").Value = 4550

# --- Row 42: Oman -------------------------------------------------------------
$ws.Range("B42").Value = 95339
$ws.Range("C42").Value = 628
$ws.Range("D42").Value = 86482
$ws.Range("E42").Value = 7982
$ws.Range("G42").Value = 10
$ws.Range("H42").Value = 875

# --- Row 47: Polonia ----------------------------------------------------------
$ws.Range("B47").Value = 81673
$ws.Range("C47").Value = 974
$ws.Range("D47").Value = 65561
$ws.Range("E47").Value = 13768
$ws.Range("G47").Value = 28
$ws.Range("H47").Value = 2344

# --- Rows 98-99: Malasia moves above Guinea (higher updated total) ----------
$ws.Range("A98").Value = "Malasia"
$ws.Range("B98").Value = 10505
$ws.Range("C98").Value = 147
$ws.Range("D98").Value = 9602
$ws.Range("E98").Value = 770
$ws.Range("G98").Value = 3
$ws.Range("H98").Value = 133

$ws.Range("A99").Value = "Guinea"
$ws.Range("B99").Value = 10387
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 9780
$ws.Range("E99").Value = 542
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 65

# --- Row 122: Hong Kong -------------------------------------------------------
$ws.Range("B122").Value = 5050
$ws.Range("C122").Value = 3
$ws.Range("D122").Value = 4749
$ws.Range("E122").Value = 198

# --- Row 142: Sri Lanka --------------------------------------------------------
$ws.Range("D142").Value = 3129
$ws.Range("E142").Value = 171

# --- Row 145: Malta -------------------------------------------------------------
$ws.Range("E145").Value = 676
$ws.Range("G145").Value = 2
$ws.Range("H145").Value = 25

# --- Row 168: Vietnam -----------------------------------------------------------
$ws.Range("D168").Value = 980
$ws.Range("E168").Value = 53

# --- Row 179: Islas Feroe ---------------------------------------------------------
$ws.Range("B179").Value = 451
$ws.Range("C179").Value = 3
$ws.Range("D179").Value = 416
$ws.Range("E179").Value = 35

# --- Rows 214-215: Montserrat moves above Islas Malvinas (tie-break) -------------
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
